$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-21 Thursday", "2023-12-22 Friday"),
    @("52×64=3328", "44×39=1716"),
    @("12×47=564", "33×69=2277"),
    @("56×90=5040", "15×79=1185"),
    @("79×58=4582", "29×58=1682"),
    @("41×88=3608", "36×47=1692"),
    @("96×41=3936", "74×79=5846"),
    @("77×36=2772", "13×40=520"),
    @("85×57=4845", "37×64=2368"),
    @("65×37=2405", "92×42=3864"),
    @("46×24=1104", "46×91=4186"),
    @("14×24=336", "80×47=3760"),
    @("88×17=1496", "98×16=1568"),
    @("46×59=2714", "45×53=2385"),
    @("74×71=5254", "56×93=5208"),
    @("98×82=8036", "26×88=2288"),
    @("79×11=869", "60×12=720"),
    @("48×92=4416", "66×55=3630"),
    @("96×48=4608", "21×33=693"),
    @("66×38=2508", "40×49=1960"),
    @("87×22=1914", "22×88=1936"),
    @("13×15=195", "46×27=1242"),
    @("21×84=1764", "62×11=682"),
    @("71×93=6603", "53×25=1325"),
    @("51×64=3264", "19×47=893"),
    @("75×57=4275", "51×36=1836")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
